# ============================================================
# khl_stats_1369_ext.xlsx refresh: 2025-11-25 -> 2025-11-26 data pull
# ============================================================
$wb = $excel.ActiveWorkbook

# --- Matches_SOG: correct F317, append new matches 322-325 ---
$wsMatches = $wb.Worksheets.Item("Matches_SOG")
$wsMatches.Range("F317").Value = 32

$wsMatches.Cells.Item(322, 1).Value = "'897818"
$wsMatches.Cells.Item(322, 2).Value = "2025-11-26T12:30:00"
$wsMatches.Cells.Item(322, 3).Value = "Адмирал"
$wsMatches.Cells.Item(322, 4).Value = "Амур"
$wsMatches.Cells.Item(322, 5).Value = 42
$wsMatches.Cells.Item(322, 6).Value = 27
$wsMatches.Cells.Item(322, 7).Value = "khl_text"

$wsMatches.Cells.Item(323, 1).Value = "'897816"
$wsMatches.Cells.Item(323, 2).Value = "2025-11-26T17:00:00"
$wsMatches.Cells.Item(323, 3).Value = "Трактор"
$wsMatches.Cells.Item(323, 4).Value = "Драконы"
$wsMatches.Cells.Item(323, 5).Value = 48
$wsMatches.Cells.Item(323, 6).Value = 27
$wsMatches.Cells.Item(323, 7).Value = "khl_text"

$wsMatches.Cells.Item(324, 1).Value = "'897817"
$wsMatches.Cells.Item(324, 2).Value = "2025-11-26T19:30:00"
$wsMatches.Cells.Item(324, 3).Value = "Северсталь"
$wsMatches.Cells.Item(324, 4).Value = "СКА"
$wsMatches.Cells.Item(324, 5).Value = 35
$wsMatches.Cells.Item(324, 6).Value = 24
$wsMatches.Cells.Item(324, 7).Value = "khl_text"

$wsMatches.Cells.Item(325, 1).Value = "'897815"
$wsMatches.Cells.Item(325, 2).Value = "2025-11-26T19:30:00"
$wsMatches.Cells.Item(325, 3).Value = "Динамо М"
$wsMatches.Cells.Item(325, 4).Value = "Локомотив"
$wsMatches.Cells.Item(325, 5).Value = 32
$wsMatches.Cells.Item(325, 6).Value = 23
$wsMatches.Cells.Item(325, 7).Value = "khl_text"

# --- Shots_HA: refresh as_of_utc + rolling totals ---
$wsShotsHA = $wb.Worksheets.Item("Shots_HA")
$wsShotsHA.Range("D2").Value = "2025-11-26T19:30:00Z"
$wsShotsHA.Range("D3").Value = "2025-11-26T19:30:00Z"
$wsShotsHA.Range("D4").Value = "2025-11-26T19:30:00Z"
$wsShotsHA.Range("E4").Value = 12
$wsShotsHA.Range("G4").Value = 456
$wsShotsHA.Range("H4").Value = 314
$wsShotsHA.Range("I4").Value = 38
$wsShotsHA.Range("J4").Value = 26.2
$wsShotsHA.Range("D5").Value = "2025-11-26T19:30:00Z"
$wsShotsHA.Range("D6").Value = "2025-11-26T19:30:00Z"
$wsShotsHA.Range("F6").Value = 16
$wsShotsHA.Range("K6").Value = 438
$wsShotsHA.Range("L6").Value = 614
$wsShotsHA.Range("M6").Value = 27.4
$wsShotsHA.Range("N6").Value = 38.4
$wsShotsHA.Range("D7").Value = "2025-11-26T19:30:00Z"
$wsShotsHA.Range("D8").Value = "2025-11-26T19:30:00Z"
$wsShotsHA.Range("E8").Value = 15
$wsShotsHA.Range("G8").Value = 499
$wsShotsHA.Range("H8").Value = 397
$wsShotsHA.Range("I8").Value = 33.3
$wsShotsHA.Range("J8").Value = 26.5
$wsShotsHA.Range("D9").Value = "2025-11-26T19:30:00Z"
$wsShotsHA.Range("D10").Value = "2025-11-26T19:30:00Z"
$wsShotsHA.Range("F10").Value = 17
$wsShotsHA.Range("K10").Value = 445
$wsShotsHA.Range("L10").Value = 618
$wsShotsHA.Range("M10").Value = 26.2
$wsShotsHA.Range("N10").Value = 36.4
$wsShotsHA.Range("D11").Value = "2025-11-26T19:30:00Z"
$wsShotsHA.Range("D12").Value = "2025-11-26T19:30:00Z"
$wsShotsHA.Range("F12").Value = 16
$wsShotsHA.Range("K12").Value = 492
$wsShotsHA.Range("L12").Value = 402
$wsShotsHA.Range("M12").Value = 30.8
$wsShotsHA.Range("N12").Value = 25.1
$wsShotsHA.Range("D13").Value = "2025-11-26T19:30:00Z"
$wsShotsHA.Range("D14").Value = "2025-11-26T19:30:00Z"
$wsShotsHA.Range("D15").Value = "2025-11-26T19:30:00Z"
$wsShotsHA.Range("F15").Value = 13
$wsShotsHA.Range("K15").Value = 384
$wsShotsHA.Range("L15").Value = 421
$wsShotsHA.Range("M15").Value = 29.5
$wsShotsHA.Range("N15").Value = 32.4
$wsShotsHA.Range("D16").Value = "2025-11-26T19:30:00Z"
$wsShotsHA.Range("D17").Value = "2025-11-26T19:30:00Z"
$wsShotsHA.Range("E17").Value = 15
$wsShotsHA.Range("G17").Value = 435
$wsShotsHA.Range("H17").Value = 343
$wsShotsHA.Range("I17").Value = 29
$wsShotsHA.Range("J17").Value = 22.9
$wsShotsHA.Range("D18").Value = "2025-11-26T19:30:00Z"
$wsShotsHA.Range("D19").Value = "2025-11-26T19:30:00Z"
$wsShotsHA.Range("D20").Value = "2025-11-26T19:30:00Z"
$wsShotsHA.Range("D21").Value = "2025-11-26T19:30:00Z"
$wsShotsHA.Range("E21").Value = 15
$wsShotsHA.Range("G21").Value = 560
$wsShotsHA.Range("H21").Value = 436
$wsShotsHA.Range("I21").Value = 37.3
$wsShotsHA.Range("J21").Value = 29.1
$wsShotsHA.Range("D22").Value = "2025-11-26T19:30:00Z"
$wsShotsHA.Range("D23").Value = "2025-11-26T19:30:00Z"

# --- Shots_Summary: refresh as_of_utc + rolling totals ---
$wsShotsSummary = $wb.Worksheets.Item("Shots_Summary")
$wsShotsSummary.Range("D2").Value = "2025-11-26T19:30:00Z"
$wsShotsSummary.Range("D3").Value = "2025-11-26T19:30:00Z"
$wsShotsSummary.Range("D4").Value = "2025-11-26T19:30:00Z"
$wsShotsSummary.Range("E4").Value = 28
$wsShotsSummary.Range("F4").Value = 953
$wsShotsSummary.Range("G4").Value = 758
$wsShotsSummary.Range("H4").Value = 34
$wsShotsSummary.Range("D5").Value = "2025-11-26T19:30:00Z"
$wsShotsSummary.Range("D6").Value = "2025-11-26T19:30:00Z"
$wsShotsSummary.Range("E6").Value = 30
$wsShotsSummary.Range("F6").Value = 857
$wsShotsSummary.Range("G6").Value = 1104
$wsShotsSummary.Range("I6").Value = 36.8
$wsShotsSummary.Range("D7").Value = "2025-11-26T19:30:00Z"
$wsShotsSummary.Range("D8").Value = "2025-11-26T19:30:00Z"
$wsShotsSummary.Range("E8").Value = 29
$wsShotsSummary.Range("F8").Value = 884
$wsShotsSummary.Range("G8").Value = 833
$wsShotsSummary.Range("H8").Value = 30.5
$wsShotsSummary.Range("I8").Value = 28.7
$wsShotsSummary.Range("D9").Value = "2025-11-26T19:30:00Z"
$wsShotsSummary.Range("D10").Value = "2025-11-26T19:30:00Z"
$wsShotsSummary.Range("E10").Value = 29
$wsShotsSummary.Range("F10").Value = 787
$wsShotsSummary.Range("G10").Value = 1034
$wsShotsSummary.Range("I10").Value = 35.7
$wsShotsSummary.Range("D11").Value = "2025-11-26T19:30:00Z"
$wsShotsSummary.Range("D12").Value = "2025-11-26T19:30:00Z"
$wsShotsSummary.Range("E12").Value = 32
$wsShotsSummary.Range("F12").Value = 1016
$wsShotsSummary.Range("G12").Value = 823
$wsShotsSummary.Range("H12").Value = 31.8
$wsShotsSummary.Range("I12").Value = 25.7
$wsShotsSummary.Range("D13").Value = "2025-11-26T19:30:00Z"
$wsShotsSummary.Range("D14").Value = "2025-11-26T19:30:00Z"
$wsShotsSummary.Range("D15").Value = "2025-11-26T19:30:00Z"
$wsShotsSummary.Range("E15").Value = 28
$wsShotsSummary.Range("F15").Value = 872
$wsShotsSummary.Range("G15").Value = 917
$wsShotsSummary.Range("H15").Value = 31.1
$wsShotsSummary.Range("I15").Value = 32.8
$wsShotsSummary.Range("D16").Value = "2025-11-26T19:30:00Z"
$wsShotsSummary.Range("D17").Value = "2025-11-26T19:30:00Z"
$wsShotsSummary.Range("E17").Value = 30
$wsShotsSummary.Range("F17").Value = 928
$wsShotsSummary.Range("G17").Value = 745
$wsShotsSummary.Range("H17").Value = 30.9
$wsShotsSummary.Range("I17").Value = 24.8
$wsShotsSummary.Range("D18").Value = "2025-11-26T19:30:00Z"
$wsShotsSummary.Range("D19").Value = "2025-11-26T19:30:00Z"
$wsShotsSummary.Range("D20").Value = "2025-11-26T19:30:00Z"
$wsShotsSummary.Range("D21").Value = "2025-11-26T19:30:00Z"
$wsShotsSummary.Range("E21").Value = 30
$wsShotsSummary.Range("F21").Value = 1070
$wsShotsSummary.Range("G21").Value = 912
$wsShotsSummary.Range("H21").Value = 35.7
$wsShotsSummary.Range("I21").Value = 30.4
$wsShotsSummary.Range("D22").Value = "2025-11-26T19:30:00Z"
$wsShotsSummary.Range("D23").Value = "2025-11-26T19:30:00Z"

# --- Meta_ext: bump as_of_utc + build_version ---
$wsMeta = $wb.Worksheets.Item("Meta_ext")
$wsMeta.Range("B2").Value = "2025-11-26T19:30:00Z"
$wsMeta.Range("D2").Value = 15
